$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

$ws.Range("G40").Value = "Warehouse employee for current user with no allowance"
$ws.Range("G61").Value = "Warehouse employee for current user with allowance"
$ws.Range("G68").Value = "Warehouse employee for current user with allowance"

$ws.Activate()
$ws.Range("G69").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
